$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the fraction-remaining input (F18) for the second table (rows 17-23)
$ws.Range("F18").Value = 0.7

# Apply the built-in "Input" cell style to the input cells of the second table
$ws.Range("C18:C22").Style = "Input"
$ws.Range("F18").Style = "Input"

# Add a total-probability check cell summing the H column
$ws.Range("H23").Formula = "=SUM(H18:H22)"
$ws.Range("H23").Style = "Normal"

# Restore the active selection to match the edited cell
$ws.Range("F19").Select()
